{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the \"Questions to answer\" list item asking about the parameters\n// that affect shares number - the new question is inserted right after it.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"What are the most important parameters that affect shares number?\") !== -1) {\n    target = paragraphs.items[i];\n  }\n}\n\nif (target) {\n  // insertParagraph after the target inherits its formatting (List\n  // Paragraph style + numbered-list numPr), matching the sibling items.\n  target.insertParagraph(\n    \"What is the best criterion for comparing models? (BIC, Cp, adjusted R^2?)\",\n    Word.InsertLocation.after\n  );\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"Questions to answer\" list item that asks about the most\n# important parameters affecting shares number - the new question about\n# the best model-comparison criterion is inserted right after it.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*What are the most important parameters that affect shares number?*\") {\n        $target = $p\n    }\n}\n\nif ($target -ne $null) {\n    # Insert a new paragraph right after the target; Word clones the\n    # paragraph/list formatting (pStyle \"a3\", numId 2 numbering, rPr lang)\n    # from $target automatically.\n    $target.Range.InsertParagraphAfter()\n    $newPara = $target.Next()\n    $newPara.Range.Text = \"What is the best criterion for comparing models? (BIC, Cp, adjusted R^2?)\"\n}\n"}
